# Swap the presentation's theme color scheme from "Integral" to the
# stock "Office Theme" palette (ppt/theme/theme1.xml: a:clrScheme).
#
# Order of ThemeColorScheme items (1-based, matches the OOXML
# <a:clrScheme> child order): dk1, lt1, dk2, lt2, accent1, accent2,
# accent3, accent4, accent5, accent6, hlink, folHlink.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

$officeColors = @(
    (RGB 0x00 0x00 0x00),  # 1  dk1
    (RGB 0xFF 0xFF 0xFF),  # 2  lt1
    (RGB 0x44 0x54 0x6A),  # 3  dk2
    (RGB 0xE7 0xE6 0xE6),  # 4  lt2
    (RGB 0x5B 0x9B 0xD5),  # 5  accent1
    (RGB 0xED 0x7D 0x31),  # 6  accent2
    (RGB 0xA5 0xA5 0xA5),  # 7  accent3
    (RGB 0xFF 0xC0 0x00),  # 8  accent4
    (RGB 0x44 0x72 0xC4),  # 9  accent5
    (RGB 0x70 0xAD 0x47),  # 10 accent6
    (RGB 0x05 0x63 0xC1),  # 11 hlink
    (RGB 0x95 0x4F 0x72)   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
